# Weekly update: insert a new price record as the first data row (row 337)
# for "Terminal La Palmera de La Serena - Papa". This pushes the existing
# rows 337-363 down to 338-364 (dimension grows from R363 to R364), and
# populates the newly inserted row 337 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 337, shifting rows 337:363 down
# to 338:364 (matching Excel's normal "Insert Cells > Entire Row" behavior).
$ws.Rows("337:337").Insert()

# Populate the newly inserted row 337 with the new record.
$ws.Cells.Item(337, 1).Value  = 8
$ws.Cells.Item(337, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(337, 3).Value  = "Coquimbo"
$ws.Cells.Item(337, 4).Value  = 44610
$ws.Cells.Item(337, 5).Value  = 4
$ws.Cells.Item(337, 6).Value  = 100114001
$ws.Cells.Item(337, 7).Value  = "Papa"
$ws.Cells.Item(337, 8).Value  = "Asterix"
$ws.Cells.Item(337, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(337, 10).Value = 2460
$ws.Cells.Item(337, 11).Value = 9000
$ws.Cells.Item(337, 12).Value = 10000
$ws.Cells.Item(337, 13).Value = 9500
$ws.Cells.Item(337, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(337, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(337, 16).Value = 380
$ws.Cells.Item(337, 17).Value = 25
$ws.Cells.Item(337, 18).Value = "Hortaliza"
